$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing grade values (5) for several students/rows
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 5

$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 5

$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 5

$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 5

$ws.Range("D19").Value = 5

$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 5

$ws.Range("F23").Value = 5

# Update the selected/active cell in the bottom-right frozen pane
$ws.Range("E19").Select()
